$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" property to the new generation time
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2025-10-30T16:36:55+00:00"

# ---------------------------------------------------------------------
# 2) Elements sheet: add a new "EntiteJuridique" reference row (row 29)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Column K ("Type(s)") now needs to be much wider to fit the new URL value
$ws.Columns.Item(11).ColumnWidth = 57.44921875

# Seed row 29 from row 28's formatting/content (keeps the shared border,
# fill, font and "empty string" cell layout used throughout the table)
$ws.Range("A28:AJ28").Copy($ws.Range("A29:AJ29"))

# Text fields for the new DispositifMedical.EntiteJuridique element
$ws.Cells.Item(29, 1).Value = "DispositifMedical.EntiteJuridique"
$ws.Cells.Item(29, 2).Value = "DispositifMedical.EntiteJuridique"
$ws.Cells.Item(29, 11).Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/EntiteJuridique`n"
$ws.Cells.Item(29, 12).Value = "Lien vers la classe EntiteJuridique"
$ws.Cells.Item(29, 13).Value = "Lien vers la classe EntiteJuridique"
$ws.Cells.Item(29, 32).Value = "DispositifMedical.EntiteJuridique"

# Min/Max (Base Min/Max) columns hold the numeric-looking text "1" — write
# them through a helper cell forced to Text via T() so Excel does not
# auto-coerce the literal into a real number (and without touching
# NumberFormat, which would otherwise mint a brand-new, unused cell style).
$helper = $ws.Cells.Item(100, 1)
$targets = @(6, 7, 33, 34)
foreach ($col in $targets) {
    $helper.Formula = '=T("1")'
    $helper.Copy()
    $ws.Cells.Item(29, $col).PasteSpecial(-4163)
}
$helper.EntireRow.Delete()

Write-Output "EntiteJuridique row + date update applied"
